$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 8.489835333333334
$ws.Range("N2").Value = 25.469506
$ws.Range("O2").Value = 0.2075776945087381
$ws.Range("P2").Value = 0.2075776945087381
$ws.Range("Q2").Value = 2.241214649976
$ws.Range("R2").Value = 20.170931849784
$ws.Range("S2").Value = 0.2075776945087381
$ws.Range("T2").Value = 0.2075776945087381

# Row 3
$ws.Range("O3").Value = 0.3214784855238645
$ws.Range("P3").Value = 0.3214784855238645
$ws.Range("S3").Value = 0.3214784855238645
$ws.Range("T3").Value = 0.3214784855238645

# Row 4
$ws.Range("M4").Value = 5.630791333333334
$ws.Range("N4").Value = 16.892374
$ws.Range("O4").Value = 0.1376736576555254
$ws.Range("P4").Value = 0.1376736576555254
$ws.Range("Q4").Value = 1.486461342504
$ws.Range("R4").Value = 13.378152082536
$ws.Range("S4").Value = 0.1376736576555254
$ws.Range("T4").Value = 0.1376736576555254

# Row 5
$ws.Range("M5").Value = 6.738585333333333
$ws.Range("N5").Value = 20.215756
$ws.Range("O5").Value = 0.1647593802263456
$ws.Range("P5").Value = 0.1647593802263456
$ws.Range("Q5").Value = 1.778905664976
$ws.Range("R5").Value = 16.010150984784
$ws.Range("S5").Value = 0.1647593802263456
$ws.Range("T5").Value = 0.1647593802263456

# Row 6
$ws.Range("M6").Value = 0.9376886666666667
$ws.Range("N6").Value = 2.813066
$ws.Range("O6").Value = 0.02292662271427322
$ws.Range("P6").Value = 0.02292662271427321
$ws.Range("Q6").Value = 0.247538555736
$ws.Range("R6").Value = 2.227847001624
$ws.Range("S6").Value = 0.02292662271427322
$ws.Range("T6").Value = 0.02292662271427321

# Row 7
$ws.Range("M7").Value = 5.954327333333333
$ws.Range("N7").Value = 17.862982
$ws.Range("O7").Value = 0.1455841593712531
$ws.Range("P7").Value = 0.1455841593712531
$ws.Range("Q7").Value = 1.571870964072
$ws.Range("R7").Value = 14.146838676648
$ws.Range("S7").Value = 0.1455841593712531
$ws.Range("T7").Value = 0.1455841593712531
